$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 0.001
$ws.Range("K9").Value = 504
$ws.Range("L9").Value = 0.00168
